$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.915.32"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.418.04"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("D5").Value = "'562.88"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").Value = "'142.85"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("E11").Value = "  -3.22%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "'25.72"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "2.854.85"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "61.824.10"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "2.411.18"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "'322.97"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "'6.82"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'66.42"
$ws.Range("E23").Value = "  +2.27%  "
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("D25").Value = "'8.81"
$ws.Range("E25").Value = "  -2.99%  "
$ws.Range("D26").Value = "'557.84"
$ws.Range("E26").Value = "  -2.77%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.532.46"
$ws.Range("E28").Value = "  -4.23%  "
$ws.Range("D29").Value = "0.0₃0934"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").Value = "'8.17"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("E34").Value = "  -3.09%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").Value = "'0.379"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").Value = "'153.72"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("E39").Value = "  -3.10%  "
$ws.Range("D40").Value = "'18.48"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").Value = "'1.82"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "'0.994"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "'147.16"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'0.0525"
$ws.Range("E46").Value = "  -2.02%  "
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").Value = "'19.76"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").Value = "'0.0920"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  +0.85%  "
